$d = $word.ActiveDocument

# --- Update the date heading (first paragraph, outside the table) ---
$dateRange = $d.Paragraphs.Item(1).Range
if ($dateRange.Text.TrimEnd([char]13,[char]7) -eq "2023-09-05 Tuesday") {
    $dateRange.Text = "2023-09-06 Wednesday"
}

# --- Update each arithmetic-problem cell in the 20x5 table ---
# Setting Range.Text directly (rather than Find.Execute) keeps each edit scoped
# exactly to its own cell/paragraph and preserves the run formatting (rFonts, sz).
$t = $d.Tables.Item(1)
$cellPara = $t.Cell(1,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "6+15=") {
    $cellPara.Range.Text = "38+59="
}
$cellPara = $t.Cell(1,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "53+39=") {
    $cellPara.Range.Text = "34-5="
}
$cellPara = $t.Cell(1,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "16+19=") {
    $cellPara.Range.Text = "6+67="
}
$cellPara = $t.Cell(1,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "17+8=") {
    $cellPara.Range.Text = "59+28="
}
$cellPara = $t.Cell(1,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "6+87=") {
    $cellPara.Range.Text = "39+23="
}
$cellPara = $t.Cell(2,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "6+85=") {
    $cellPara.Range.Text = "84+7="
}
$cellPara = $t.Cell(2,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "52-36=") {
    $cellPara.Range.Text = "38+28="
}
$cellPara = $t.Cell(2,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "23-18=") {
    $cellPara.Range.Text = "90-17="
}
$cellPara = $t.Cell(2,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "12-6=") {
    $cellPara.Range.Text = "90-83="
}
$cellPara = $t.Cell(2,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "29+59=") {
    $cellPara.Range.Text = "90-16="
}
$cellPara = $t.Cell(3,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "23+59=") {
    $cellPara.Range.Text = "43+49="
}
$cellPara = $t.Cell(3,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "72-27=") {
    $cellPara.Range.Text = "46+15="
}
$cellPara = $t.Cell(3,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "76-8=") {
    $cellPara.Range.Text = "64+8="
}
$cellPara = $t.Cell(3,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "73-69=") {
    $cellPara.Range.Text = "54-9="
}
$cellPara = $t.Cell(3,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "82-28=") {
    $cellPara.Range.Text = "82-67="
}
$cellPara = $t.Cell(4,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "4+8=") {
    $cellPara.Range.Text = "81-72="
}
$cellPara = $t.Cell(4,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "66-7=") {
    $cellPara.Range.Text = "92-33="
}
$cellPara = $t.Cell(4,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "61-49=") {
    $cellPara.Range.Text = "30-25="
}
$cellPara = $t.Cell(4,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "39+28=") {
    $cellPara.Range.Text = "8+88="
}
$cellPara = $t.Cell(4,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "60-51=") {
    $cellPara.Range.Text = "58+17="
}
$cellPara = $t.Cell(5,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "37-18=") {
    $cellPara.Range.Text = "94-9="
}
$cellPara = $t.Cell(5,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "51-19=") {
    $cellPara.Range.Text = "36+29="
}
$cellPara = $t.Cell(5,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "51-12=") {
    $cellPara.Range.Text = "72-56="
}
$cellPara = $t.Cell(5,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "21-18=") {
    $cellPara.Range.Text = "60-16="
}
$cellPara = $t.Cell(5,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "91-74=") {
    $cellPara.Range.Text = "25-8="
}
$cellPara = $t.Cell(6,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "4+48=") {
    $cellPara.Range.Text = "8+48="
}
$cellPara = $t.Cell(6,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "44-6=") {
    $cellPara.Range.Text = "83+9="
}
$cellPara = $t.Cell(6,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "42-27=") {
    $cellPara.Range.Text = "45+16="
}
$cellPara = $t.Cell(6,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "32+29=") {
    $cellPara.Range.Text = "70-32="
}
$cellPara = $t.Cell(6,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "60-31=") {
    $cellPara.Range.Text = "8+76="
}
$cellPara = $t.Cell(7,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "74-19=") {
    $cellPara.Range.Text = "81-44="
}
$cellPara = $t.Cell(7,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "24+8=") {
    $cellPara.Range.Text = "71-57="
}
$cellPara = $t.Cell(7,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "6+6=") {
    $cellPara.Range.Text = "77-8="
}
$cellPara = $t.Cell(7,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "27+17=") {
    $cellPara.Range.Text = "73-35="
}
$cellPara = $t.Cell(7,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "96-58=") {
    $cellPara.Range.Text = "5+86="
}
$cellPara = $t.Cell(8,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "34+19=") {
    $cellPara.Range.Text = "82-9="
}
$cellPara = $t.Cell(8,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "94-16=") {
    $cellPara.Range.Text = "65+16="
}
$cellPara = $t.Cell(8,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "71-6=") {
    $cellPara.Range.Text = "53+19="
}
$cellPara = $t.Cell(8,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "52-27=") {
    $cellPara.Range.Text = "59+7="
}
$cellPara = $t.Cell(8,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "73-16=") {
    $cellPara.Range.Text = "78-59="
}
$cellPara = $t.Cell(9,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "9+75=") {
    $cellPara.Range.Text = "87+7="
}
$cellPara = $t.Cell(9,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "2+39=") {
    $cellPara.Range.Text = "77+6="
}
$cellPara = $t.Cell(9,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "66+28=") {
    $cellPara.Range.Text = "48+5="
}
$cellPara = $t.Cell(9,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "63-59=") {
    $cellPara.Range.Text = "84-16="
}
$cellPara = $t.Cell(9,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "45+27=") {
    $cellPara.Range.Text = "18+63="
}
$cellPara = $t.Cell(10,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "51-37=") {
    $cellPara.Range.Text = "38+46="
}
$cellPara = $t.Cell(10,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "23+9=") {
    $cellPara.Range.Text = "16-9="
}
$cellPara = $t.Cell(10,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "33+59=") {
    $cellPara.Range.Text = "92-17="
}
$cellPara = $t.Cell(10,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "86-8=") {
    $cellPara.Range.Text = "37+49="
}
$cellPara = $t.Cell(10,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "80-14=") {
    $cellPara.Range.Text = "4+39="
}
$cellPara = $t.Cell(11,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "73-46=") {
    $cellPara.Range.Text = "41-12="
}
$cellPara = $t.Cell(11,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "65+29=") {
    $cellPara.Range.Text = "66+9="
}
$cellPara = $t.Cell(11,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "82-47=") {
    $cellPara.Range.Text = "37+44="
}
$cellPara = $t.Cell(11,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "7+76=") {
    $cellPara.Range.Text = "55-19="
}
$cellPara = $t.Cell(11,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "45+17=") {
    $cellPara.Range.Text = "19+34="
}
$cellPara = $t.Cell(12,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "36+17=") {
    $cellPara.Range.Text = "41-6="
}
$cellPara = $t.Cell(12,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "73-9=") {
    $cellPara.Range.Text = "34-17="
}
$cellPara = $t.Cell(12,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "28+63=") {
    $cellPara.Range.Text = "59+18="
}
$cellPara = $t.Cell(12,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "36+15=") {
    $cellPara.Range.Text = "55-38="
}
$cellPara = $t.Cell(12,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "95-38=") {
    $cellPara.Range.Text = "83-28="
}
$cellPara = $t.Cell(13,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "33-28=") {
    $cellPara.Range.Text = "66-59="
}
$cellPara = $t.Cell(13,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "8+4=") {
    $cellPara.Range.Text = "68-39="
}
$cellPara = $t.Cell(13,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "20-19=") {
    $cellPara.Range.Text = "5+86="
}
$cellPara = $t.Cell(13,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "48+13=") {
    $cellPara.Range.Text = "50-7="
}
$cellPara = $t.Cell(13,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "18+36=") {
    $cellPara.Range.Text = "16+5="
}
$cellPara = $t.Cell(14,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "96-89=") {
    $cellPara.Range.Text = "31-27="
}
$cellPara = $t.Cell(14,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "48+33=") {
    $cellPara.Range.Text = "54+17="
}
$cellPara = $t.Cell(14,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "47+5=") {
    $cellPara.Range.Text = "66-59="
}
$cellPara = $t.Cell(14,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "16+16=") {
    $cellPara.Range.Text = "19+66="
}
$cellPara = $t.Cell(14,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "17+19=") {
    $cellPara.Range.Text = "37+36="
}
$cellPara = $t.Cell(15,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "34+9=") {
    $cellPara.Range.Text = "46+19="
}
$cellPara = $t.Cell(15,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "27+7=") {
    $cellPara.Range.Text = "70-58="
}
$cellPara = $t.Cell(15,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "28+28=") {
    $cellPara.Range.Text = "21-7="
}
$cellPara = $t.Cell(15,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "23+48=") {
    $cellPara.Range.Text = "70-48="
}
$cellPara = $t.Cell(15,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "92-43=") {
    $cellPara.Range.Text = "93-76="
}
$cellPara = $t.Cell(16,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "7+64=") {
    $cellPara.Range.Text = "93-14="
}
$cellPara = $t.Cell(16,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "17+56=") {
    $cellPara.Range.Text = "12+59="
}
$cellPara = $t.Cell(16,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "90-74=") {
    $cellPara.Range.Text = "6+57="
}
$cellPara = $t.Cell(16,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "92-28=") {
    $cellPara.Range.Text = "29+46="
}
$cellPara = $t.Cell(16,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "83-26=") {
    $cellPara.Range.Text = "9+59="
}
$cellPara = $t.Cell(17,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "17+67=") {
    $cellPara.Range.Text = "72-14="
}
$cellPara = $t.Cell(17,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "73+8=") {
    $cellPara.Range.Text = "23+38="
}
$cellPara = $t.Cell(17,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "42-24=") {
    $cellPara.Range.Text = "53-7="
}
$cellPara = $t.Cell(17,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "52-33=") {
    $cellPara.Range.Text = "93-14="
}
$cellPara = $t.Cell(17,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "81-36=") {
    $cellPara.Range.Text = "72-4="
}
$cellPara = $t.Cell(18,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "19+55=") {
    $cellPara.Range.Text = "38+59="
}
$cellPara = $t.Cell(18,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "80-68=") {
    $cellPara.Range.Text = "31-8="
}
$cellPara = $t.Cell(18,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "31-22=") {
    $cellPara.Range.Text = "18+54="
}
$cellPara = $t.Cell(18,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "82-26=") {
    $cellPara.Range.Text = "4+28="
}
$cellPara = $t.Cell(18,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "4+77=") {
    $cellPara.Range.Text = "9+59="
}
$cellPara = $t.Cell(19,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "78-39=") {
    $cellPara.Range.Text = "95-87="
}
$cellPara = $t.Cell(19,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "38+26=") {
    $cellPara.Range.Text = "61-54="
}
$cellPara = $t.Cell(19,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "32-27=") {
    $cellPara.Range.Text = "74-48="
}
$cellPara = $t.Cell(19,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "9+32=") {
    $cellPara.Range.Text = "64-35="
}
$cellPara = $t.Cell(19,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "80-19=") {
    $cellPara.Range.Text = "27+48="
}
$cellPara = $t.Cell(20,1).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "84-68=") {
    $cellPara.Range.Text = "75-9="
}
$cellPara = $t.Cell(20,2).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "33-5=") {
    $cellPara.Range.Text = "66-57="
}
$cellPara = $t.Cell(20,3).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "68+28=") {
    $cellPara.Range.Text = "9+85="
}
$cellPara = $t.Cell(20,4).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "47+14=") {
    $cellPara.Range.Text = "79+15="
}
$cellPara = $t.Cell(20,5).Range.Paragraphs.Item(1)
if ($cellPara.Range.Text.TrimEnd([char]13,[char]7) -eq "70-26=") {
    $cellPara.Range.Text = "15+56="
}
